# #CRM-177 District Field in SF list Excel
#
# Adds a new "District" column to the vendor export template, positioned
# right after the existing "ID" / Appliances column (new column N), shifting
# every column from the old N onward one position to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at N (14th column), shifting M..AH -> N..AI
$ws.Range("N1").EntireColumn.Insert()

# Match the width of the new column to its left-hand neighbour (column M)
# as closely as this runtime's ColumnWidth precision allows.
$ws.Range("N1").ColumnWidth = $ws.Range("M1").ColumnWidth

# Populate the new header / placeholder cells
$ws.Range("N1").Value = "District"
$ws.Range("N2").Value = "{vendor:district}"

# Restore the selection shown in the saved file
$ws.Range("M14").Select()
